$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "P2" "320018702600"
Set-TextValue "Q2" "`$20.36"
$ws.Range("R2").Value = "FAIL"

# Row 3
Set-TextValue "P3" "320018702611"
Set-TextValue "Q3" "`$19.30"
$ws.Range("R3").Value = "FAIL"

# Row 4
Set-TextValue "P4" "320018702644"
Set-TextValue "Q4" "`$53.14"
$ws.Range("R4").Value = "FAIL"

# Row 5
Set-TextValue "P5" "320018702655"
Set-TextValue "Q5" "`$46.27"
$ws.Range("R5").Value = "FAIL"
